# The sheet currently has an extra leading column A (a plain row-index
# column holding 2 / 9) that duplicates the values already present in the
# last column F (2 / 2 and 9 / 9). The edit removes that redundant column A,
# shifting B:F left into A:E, and fixes a header label typo along the way
# (MODEL_CONDITION -> MODELCONDITION).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text before the shift (E1 here is "MODEL_CONDITION",
# which ends up as D1 once column A is removed).
$ws.Range("E1").Value = "MODELCONDITION"

# Remove column A; B:F (now holding the real data) shift left to A:E.
$ws.Columns("A").Delete()

Write-Output "Column A removed; header label corrected."
